$d = $word.ActiveDocument
$d.Content.Find.Execute("demonstration", $true, $false, $false, $false, $false,
                         $true, 1, $false, "demonstration", 2)
